$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 197; this shifts all existing rows
# 197-256 down to 198-257 (matching the target dimension A1:R257).
$ws.Rows(197).Insert()

# Populate the newly inserted row 197 with the new record.
$ws.Range("A197").Value = 10
$ws.Range("B197").Value = "Vega Modelo de Temuco"
$ws.Range("C197").Value = "La Araucanía"
$ws.Range("D197").Value = 44524
$ws.Range("E197").Value = 9
$ws.Range("F197").Value = 100112037
$ws.Range("G197").Value = "Cebollín"
$ws.Range("H197").Value = "Sin especificar"
$ws.Range("I197").Value = "Primera"
$ws.Range("J197").Value = 40
$ws.Range("K197").Value = 5000
$ws.Range("L197").Value = 5000
$ws.Range("M197").Value = 5000
$ws.Range("N197").Value = "$/docena de paquetes"
$ws.Range("O197").Value = "Región de O'Higgins"
$ws.Range("P197").Value = 417
$ws.Range("Q197").Value = 12
$ws.Range("R197").Value = "Hortaliza"
